$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings formatted like the coinranking.com source (dot-group
# separators, fixed decimals, etc.), stored as text in the workbook. Values that look
# like plain numbers ("392.57", "1.00", "0.0000105", ...) would otherwise be silently
# re-interpreted as numeric by Excel, so we force Text format for the write and then
# restore the default "Normal" style so no stray cell formatting is left behind.
$ws.Range("D2").Value = "56.885.43"
$ws.Range("E2").Value = "  +6.57%  "
$ws.Range("D3").Value = "3.223.97"
$ws.Range("E3").Value = "  +2.15%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "392.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.54%  "
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.570"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.78%  "
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "3.216.87"
$ws.Range("E8").Value = "  +2.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.611"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "38.77"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0962"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +10.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.141"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.43%  "
$ws.Range("D14").Value = "3.752.75"
$ws.Range("E14").Value = "  +3.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.05"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.08%  "
$ws.Range("D17").Value = "3.226.73"
$ws.Range("E17").Value = "  +2.10%  "
$ws.Range("E18").Value = "  -2.86%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.59%  "
$ws.Range("D20").Value = "56.798.73"
$ws.Range("E20").Value = "  +6.57%  "
$ws.Range("E21").Value = "  +1.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0000105"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "294.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.85%  "
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.168"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.74%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.35%  "
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("E32").Value = "  -2.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "10.99"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "37.14"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0481"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.17%  "
$ws.Range("E36").Value = "  +1.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.56"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.52"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "134.36"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.74%  "
$ws.Range("E42").Value = "  +1.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.94"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.81"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.280"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.76"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.91%  "
$ws.Range("D48").Value = "2.143.07"
$ws.Range("E48").Value = "  +2.52%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.10"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.20%  "
$ws.Range("E50").Value = "  +23.58%  "
$ws.Range("E51").Value = "  -2.34%  "
